$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 41, shifting rows 41:52 down to 42:53
$ws.Rows.Item(41).Insert()

# Copy the (now shifted) row 42 formatting/values into the new row 41
$ws.Rows.Item(42).Copy()
$ws.Rows.Item(41).PasteSpecial()
$excel.CutCopyMode = $false

# Update the new row's changed values (the date and price/volume figures)
$ws.Range("D41").Value = 45027
$ws.Range("J41").Value = 70
$ws.Range("K41").Value = 6000
$ws.Range("L41").Value = 6000
$ws.Range("M41").Value = 6000
$ws.Range("P41").Value = 375
